# This script applies a weekly data update: a new price-report row is
# inserted at row 14 (the most recent date), pushing all existing data
# rows down by one (old row 14 becomes row 15, ..., old row 123 becomes
# row 124). The new row 14 reuses the same market/category metadata as
# the row it displaces, with updated date and price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14; this shifts rows 14:123 down to 15:124
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new weekly record
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 45035
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112030
$ws.Range("G14").Value = "Poroto granado"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 27000
$ws.Range("L14").Value = 28000
$ws.Range("M14").Value = 27500
$ws.Range("N14").Value = "$/malla 25 kilos"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 1100
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"

# Ensure the date cell keeps the same date-number style as the rest of column D
$ws.Range("D14").NumberFormat = $ws.Range("D15").NumberFormat
